$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 190.28572
$ws.Range("I12").Value = 196.4
$ws.Range("J12").Value = 175
$ws.Range("K12").Value = 196.4
$ws.Range("L12").Value = 175
$ws.Range("M12").Value = -26.40000000000001
$ws.Range("N12").Value = -515
$ws.Range("H32").Value = 7001
$ws.Range("I32").Value = 7001
$ws.Range("K32").Value = 7001
$ws.Range("M32").Value = -6675
$ws.Range("H46").Value = 2612
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 2882.6667
$ws.Range("K46").Value = 5400
$ws.Range("L46").Value = 8648.000100000001
$ws.Range("M46").Value = -5281
$ws.Range("N46").Value = -8886.000100000001
$ws.Range("H60").Value = 2612
$ws.Range("I60").Value = 1800
$ws.Range("J60").Value = 2882.6667
$ws.Range("K60").Value = 5400
$ws.Range("L60").Value = 8648.000100000001
$ws.Range("M60").Value = -4916
$ws.Range("N60").Value = -9616.000100000001
$ws.Range("H70").Value = 893.2
$ws.Range("J70").Value = 916.5
$ws.Range("L70").Value = 2749.5
$ws.Range("N70").Value = -3289.5
$ws.Range("H73").Value = 893.2
$ws.Range("J73").Value = 916.5
$ws.Range("L73").Value = 2749.5
$ws.Range("N73").Value = -4621.5
$ws.Range("H112").Value = 1506.0322
$ws.Range("I112").Value = 1383.75
$ws.Range("J112").Value = 1524.1482
$ws.Range("K112").Value = 4151.25
$ws.Range("L112").Value = 4572.444600000001
$ws.Range("M112").Value = -3043.25
$ws.Range("N112").Value = -6788.444600000001
$ws.Range("H129").Value = 1103.0588
$ws.Range("I129").Value = 659.4545000000001
$ws.Range("K129").Value = 1978.3635
$ws.Range("M129").Value = 3021.6365
$ws.Range("H137").Value = 68445.53
$ws.Range("I137").Value = 1965.3334
$ws.Range("K137").Value = 5896.0002
$ws.Range("M137").Value = -3346.0002
$ws.Range("H138").Value = 2279.2
$ws.Range("J138").Value = 2519.25
$ws.Range("L138").Value = 7557.75
$ws.Range("N138").Value = -17837.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1395
$ws.Range("J3").Value = 1395
$ws.Range("L3").Value = 1395
$ws.Range("N3").Value = -1625
$ws.Range("H4").Value = 143359
$ws.Range("I4").Value = 154371.23
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 154371.23
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -154255.23
$ws.Range("N4").Value = -432
$ws.Range("H32").Value = 36590.965
$ws.Range("I32").Value = 21151.56
$ws.Range("K32").Value = 21151.56
$ws.Range("M32").Value = -20864.56
$ws.Range("H74").Value = 1677.5
$ws.Range("I74").Value = 1737
$ws.Range("J74").Value = 1499
$ws.Range("K74").Value = 1737
$ws.Range("L74").Value = 1499
$ws.Range("M74").Value = -863
$ws.Range("N74").Value = -3247
$ws.Range("H77").Value = 1677.5
$ws.Range("I77").Value = 1737
$ws.Range("J77").Value = 1499
$ws.Range("K77").Value = 8685
$ws.Range("L77").Value = 7495
$ws.Range("M77").Value = -4317
$ws.Range("N77").Value = -16231
$ws.Range("H122").Value = 20983.727
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 282
$ws.Range("I22").Value = 258.33334
$ws.Range("J22").Value = 329.33334
$ws.Range("K22").Value = 258.33334
$ws.Range("L22").Value = 329.33334
$ws.Range("M22").Value = 91.66665999999998
$ws.Range("N22").Value = -1029.33334
$ws.Range("H31").Value = 1242.2858
$ws.Range("I31").Value = 1242.2858
$ws.Range("K31").Value = 1242.2858
$ws.Range("M31").Value = -947.2858000000001
$ws.Range("H34").Value = 1242.2858
$ws.Range("I34").Value = 1242.2858
$ws.Range("K34").Value = 1242.2858
$ws.Range("M34").Value = -1040.2858
$ws.Range("H59").Value = 12674.889
$ws.Range("J59").Value = 13567.143
$ws.Range("L59").Value = 13567.143
$ws.Range("N59").Value = -15857.143
$ws.Range("H132").Value = 2198.1538
$ws.Range("I132").Value = 2225.3044
$ws.Range("J132").Value = 1990
$ws.Range("K132").Value = 6675.9132
$ws.Range("L132").Value = 5970
$ws.Range("M132").Value = -4145.9132
$ws.Range("N132").Value = -11030
$ws.Range("H134").Value = 1945.0294
$ws.Range("I134").Value = 1844.5555
$ws.Range("K134").Value = 5533.666499999999
$ws.Range("M134").Value = -2998.666499999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6738266.5
$ws.Range("I4").Value = 1253393.2
$ws.Range("J4").Value = 56102130
$ws.Range("K4").Value = 3760179.6
$ws.Range("L4").Value = 168306390
$ws.Range("M4").Value = -3760067.6
$ws.Range("N4").Value = -168306614
$ws.Range("H60").Value = 1499.8334
$ws.Range("I60").Value = 1999
$ws.Range("K60").Value = 5997
$ws.Range("M60").Value = -5746
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents() | Out-Null
$ws.Range("H114").Value = 28573150
$ws.Range("J114").Value = 2000
$ws.Range("L114").Value = 6000
$ws.Range("N114").Value = -12508
$ws.Range("H122").Value = 848.8
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 848.8
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7639.2
$ws.Range("M122").ClearContents() | Out-Null
$ws.Range("N122").Value = -12539.2
$ws.Range("H126").Value = 8994.75
$ws.Range("H131").Value = 2189.4707
$ws.Range("I131").Value = 1752.625
$ws.Range("J131").Value = 2577.7778
$ws.Range("K131").Value = 5257.875
$ws.Range("L131").Value = 7733.3334
$ws.Range("M131").Value = -217.875
$ws.Range("N131").Value = -17813.3334
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 6042
$ws.Range("J26").Value = 6042
$ws.Range("L26").Value = 6042
$ws.Range("N26").Value = -6602
$ws.Range("H50").Value = 6042
$ws.Range("J50").Value = 6042
$ws.Range("L50").Value = 6042
$ws.Range("N50").Value = -7038
$ws.Range("H52").Value = 40120.332
$ws.Range("J52").Value = 40120.332
$ws.Range("L52").Value = 40120.332
$ws.Range("N52").Value = -40638.332
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents() | Out-Null
$ws.Range("H58").Value = 19998.5
$ws.Range("I58").Value = 15000
$ws.Range("J58").Value = 24997
$ws.Range("K58").Value = 15000
$ws.Range("L58").Value = 24997
$ws.Range("M58").Value = -14723
$ws.Range("N58").Value = -25551
$ws.Range("H99").Value = 12057
$ws.Range("I99").Value = 9419.571
$ws.Range("K99").Value = 9419.571
$ws.Range("M99").Value = -7173.571
$ws.Range("H136").Value = 53081.25
$ws.Range("J136").Value = 53081.25
$ws.Range("L136").Value = 159243.75
$ws.Range("N136").Value = -164343.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 501.63635
$ws.Range("I16").Value = 569
$ws.Range("J16").Value = 198.5
$ws.Range("K16").Value = 569
$ws.Range("L16").Value = 198.5
$ws.Range("M16").Value = -399
$ws.Range("N16").Value = -538.5
$ws.Range("H22").Value = 1626.05
$ws.Range("I22").Value = 1613.1538
$ws.Range("J22").Value = 1650
$ws.Range("K22").Value = 1613.1538
$ws.Range("L22").Value = 1650
$ws.Range("M22").Value = -1318.1538
$ws.Range("N22").Value = -2240
$ws.Range("H27").Value = 1626.05
$ws.Range("I27").Value = 1613.1538
$ws.Range("J27").Value = 1650
$ws.Range("K27").Value = 1613.1538
$ws.Range("L27").Value = 1650
$ws.Range("M27").Value = -1506.1538
$ws.Range("N27").Value = -1864
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents() | Out-Null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents() | Out-Null
$ws.Range("H93").Value = 1505.1052
$ws.Range("I93").Value = 1485.4286
$ws.Range("J93").Value = 1560.2
$ws.Range("K93").Value = 1485.4286
$ws.Range("L93").Value = 1560.2
$ws.Range("M93").Value = -237.4286
$ws.Range("N93").Value = -4056.2
$ws.Range("H136").Value = 3462.35
$ws.Range("I136").Value = 2800.4285
$ws.Range("J136").Value = 3818.7693
$ws.Range("K136").Value = 8401.2855
$ws.Range("L136").Value = 11456.3079
$ws.Range("M136").Value = -5851.2855
$ws.Range("N136").Value = -16556.3079
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2998.3333
$ws.Range("J4").Value = 2998.3333
$ws.Range("L4").Value = 2998.3333
$ws.Range("N4").Value = -3224.3333
$ws.Range("H61").Value = 30497
$ws.Range("J61").Value = 30497
$ws.Range("L61").Value = 30497
$ws.Range("N61").Value = -31081
$ws.Range("H96").Value = 2248.625
$ws.Range("I96").Value = 1832
$ws.Range("K96").Value = 1832
$ws.Range("M96").Value = -459
